$d = $word.ActiveDocument

# 1. Append "f" to the organizationName placeholder
$d.Content.Find.Execute(
    "{d.parcels[i].owners[i].organizationName}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{d.parcels[i].owners[i].organizationName}f", 2)

# 2. Insert "convCRLF:" before "ifEM():show(.noData)}" for the specific fields
$fields = @(
    "parcelsAgricultureDescription",
    "parcelsAgricultureImprovementDescription",
    "parcelsNonAgricultureUseDescription",
    "purpose",
    "soilFillTypeToPlace",
    "soilTypeRemoved",
    "soilStructureFarmUseReason",
    "soilStructureResidentialUseReason",
    "soilAgriParcelActivity",
    "soilStructureResidentialAccessoryUseReason",
    "soilStructureOtherUseReason"
)

foreach ($field in $fields) {
    $old = "{d.$field`:ifEM():show(.noData)}"
    $new = "{d.$field`:convCRLF:ifEM():show(.noData)}"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
